$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.296.07"
$ws.Range("E2").Value = "  +0.47%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.643.54"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.39"
$ws.Range("E5").Value = "  +0.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.70"
$ws.Range("E6").Value = "  +0.52%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -0.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.642.16"
$ws.Range("E9").Value = "  +0.46%  "

$ws.Range("E10").Value = "  +8.67%  "

$ws.Range("E11").Value = "  -0.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  +0.76%  "

$ws.Range("E13").Value = "  +1.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000193"
$ws.Range("E14").Value = "  +2.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.92"
$ws.Range("E15").Value = "  +0.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.125.12"
$ws.Range("E16").Value = "  +0.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.040.14"
$ws.Range("E17").Value = "  +0.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.651.24"
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.37"
$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "364.17"
$ws.Range("E20").Value = "  -1.25%  "

$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.36"
$ws.Range("E22").Value = "  +2.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.86"
$ws.Range("E23").Value = "  +1.16%  "

$ws.Range("E24").Value = "  -1.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "75.35"
$ws.Range("E25").Value = "  +4.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.77"
$ws.Range("E27").Value = "  -0.76%  "

$ws.Range("E28").Value = "  +1.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.01"
$ws.Range("E30").Value = "  +1.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "564.83"
$ws.Range("E31").Value = "  -1.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.06"
$ws.Range("E32").Value = "  +1.66%  "

$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("E34").Value = "  +1.02%  "

$ws.Range("E35").Value = "  +2.27%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.58"
$ws.Range("E37").Value = "  +2.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.66"
$ws.Range("E38").Value = "  +1.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.33"
$ws.Range("E39").Value = "  +0.88%  "

$ws.Range("E40").Value = "  +1.61%  "

$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("E42").Value = "  -0.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₆0342"
$ws.Range("E43").Value = "  +2.29%  "

$ws.Range("E44").Value = "  -0.84%  "

$ws.Range("E45").Value = "  +2.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.64"
$ws.Range("E46").Value = "  +1.37%  "

$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.66"
$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("E49").Value = "  +1.82%  "

$ws.Range("E50").Value = "  +0.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.79"
$ws.Range("E51").Value = "  -0.91%  "
